$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Range("N3").Value = 2.1
$ws.Range("O3").Value = 1.73
# Row 4
$ws.Range("H4").Value = 4.1
$ws.Range("I4").Value = 6
$ws.Range("AE4").Value = 15
$ws.Range("AF4").Value = 29
# Row 5
$ws.Range("G5").Value = 2.45
$ws.Range("H5").Value = 2.75
$ws.Range("I5").Value = 3.3
$ws.Range("K5").Value = 4.75
$ws.Range("W5").Value = 26
$ws.Range("AE5").Value = 6
# Row 7
$ws.Range("G7").Value = 2.35
$ws.Range("H7").Value = 3.1
$ws.Range("I7").Value = 3.2
$ws.Range("T7").Value = 6.5
$ws.Range("U7").Value = 10
$ws.Range("V7").Value = 10
$ws.Range("W7").Value = 21
$ws.Range("X7").Value = 23
$ws.Range("Z7").Value = 7
$ws.Range("AA7").Value = 6
$ws.Range("AB7").Value = 17
$ws.Range("AG7").Value = 12
$ws.Range("AH7").Value = 34
$ws.Range("AI7").Value = 29
# Row 8
$ws.Range("G8").Value = 1.42
$ws.Range("I8").Value = 7
$ws.Range("W8").Value = 9.5
$ws.Range("Z8").Value = 13
$ws.Range("AD8").Value = 301
$ws.Range("AG8").Value = 21
$ws.Range("AH8").Value = 81
# Row 10
$ws.Range("K10").Value = 19
$ws.Range("N10").Value = 1.44
$ws.Range("O10").Value = 2.7
# Row 11
$ws.Range("G11").Value = 1.5
$ws.Range("H11").Value = 4.75
$ws.Range("I11").Value = 5.5
$ws.Range("L11").Value = 1.17
$ws.Range("M11").Value = 5
$ws.Range("N11").Value = 1.57
$ws.Range("O11").Value = 2.35
$ws.Range("P11").Value = 1.29
$ws.Range("Q11").Value = 3.5
$ws.Range("U11").Value = 8.5
$ws.Range("W11").Value = 11
$ws.Range("X11").Value = 11
$ws.Range("Z11").Value = 17
$ws.Range("AA11").Value = 9
$ws.Range("AB11").Value = 17
# Row 12
$ws.Range("G12").Value = 3
$ws.Range("I12").Value = 2.15
$ws.Range("L12").Value = 1.14
$ws.Range("M12").Value = 5.5
$ws.Range("U12").Value = 19
$ws.Range("X12").Value = 21
$ws.Range("AG12").Value = 9.5
$ws.Range("AJ12").Value = 21
# Row 15
$ws.Range("N15").Value = 1.6
$ws.Range("O15").Value = 2.3
# Row 18
$ws.Range("N18").Value = 1.53
$ws.Range("O18").Value = 2.4
# Row 19
$ws.Range("L19").Value = 1.18
$ws.Range("M19").Value = 4.5
$ws.Range("N19").Value = 1.6
$ws.Range("O19").Value = 2.3
$ws.Range("AE19").Value = 15
# Row 20
$ws.Range("N20").Value = 2.08
$ws.Range("O20").Value = 1.73
# Row 23
$ws.Range("G23").Value = 2.88
$ws.Range("I23").Value = 2.5
$ws.Range("W23").Value = 29
$ws.Range("AE23").Value = 7
$ws.Range("AG23").Value = 11
$ws.Range("AH23").Value = 26
# Row 28
$ws.Range("G28").Value = 2.6
$ws.Range("H28").Value = 3.2
$ws.Range("K28").Value = 7.5
$ws.Range("L28").Value = 1.44
$ws.Range("M28").Value = 2.63
$ws.Range("N28").Value = 2.35
$ws.Range("O28").Value = 1.57
$ws.Range("P28").Value = 1.53
$ws.Range("Q28").Value = 2.38
$ws.Range("R28").Value = 2
$ws.Range("S28").Value = 1.75
$ws.Range("Y28").Value = 41
$ws.Range("Z28").Value = 7.5
$ws.Range("AC28").Value = 67
$ws.Range("AD28").Value = 451
$ws.Range("AE28").Value = 7
$ws.Range("AF28").Value = 12
$ws.Range("AI28").Value = 26
$ws.Range("AJ28").Value = 41
# Row 32
$ws.Range("G32").Value = 1.7
$ws.Range("I32").Value = 4.2
$ws.Range("T32").Value = 8.25
$ws.Range("U32").Value = 8.75
$ws.Range("W32").Value = 13.5
$ws.Range("AA32").Value = 7.7
$ws.Range("AB32").Value = 14.5
$ws.Range("AE32").Value = 14
$ws.Range("AF32").Value = 25
# Row 34
$ws.Range("G34").Value = 2.7
$ws.Range("I34").Value = 2.55
$ws.Range("U34").Value = 13
$ws.Range("V34").Value = 11
# Row 35
$ws.Range("L35").Value = 1.57
$ws.Range("M35").Value = 2.25
$ws.Range("N35").Value = 2.88
$ws.Range("O35").Value = 1.4
$ws.Range("V35").Value = 12
$ws.Range("AG35").Value = 12
# Row 36
$ws.Range("G36").Value = 1.83
$ws.Range("H36").Value = 3
$ws.Range("I36").Value = 4.75
$ws.Range("N36").Value = 2.5
$ws.Range("O36").Value = 1.5
$ws.Range("P36").Value = 1.53
$ws.Range("Q36").Value = 2.38
$ws.Range("V36").Value = 9.5
$ws.Range("W36").Value = 15
$ws.Range("Y36").Value = 41
$ws.Range("Z36").Value = 6.5
$ws.Range("AA36").Value = 6
